# Changed path to data
# Update the time-bucket / count table: shift the bucket boundaries and
# counts, and append two more rows (07 and 08) that extend the table down
# to row 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 20

# Row 3
$ws.Range("A3").Value = 0.27083333333333331
$ws.Range("B3").Value = 0.375
$ws.Range("C3").Value = 400

# Row 4
$ws.Range("A4").Value = 0.375
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 100

# Row 5
$ws.Range("A5").Value = 0.5
$ws.Range("B5").Value = 0.66666666666666663
$ws.Range("C5").Value = 120

# Row 6
$ws.Range("A6").Value = 0.66666666666666663
$ws.Range("B6").Value = 0.75
$ws.Range("C6").Value = 200

# Row 7 (new) - fill-down continuation of the table, so it inherits the
# same formatting (style 6/3) as the rows above rather than the bare
# column-default style.
$ws.Range("A7").Value = 0.75
$ws.Range("B7").Value = 0.91666666666666663
$ws.Range("C7").Value = 250
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)   # xlPasteFormats

# Row 8 (new) - a fresh value typed below the table, so it just takes on
# the plain column-default style (4 for A/B, 5 for C).
$ws.Range("A8").Value = 0.91666666666666663
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 100

# Match the saved selection in the target file
$ws.Range("B8").Select() | Out-Null
